$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A values: rows 2-24 from 2 -> 1, rows 25-31 from 3 -> 2
$ws.Range("A2:A24").Value = 1
$ws.Range("A25:A31").Value = 2

# Update the active selection to B24
$ws.Range("B24").Select()
